# Weekly update: insert two new rows (new date 44841) at the top of the
# detail block (right after row 400), shifting all existing data rows
# down by two. This mirrors how the source data is produced: new rows
# for a Primera/Segunda pair are appended after the previous reporting
# period's pair, which here was inserted at the former row 401 position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the existing row 401 (old rows 401..506
# shift to 403..508).
$ws.Rows("401:402").Insert()

# --- New row 401 ("Primera") ---
$ws.Range("A401").Value2 = 8
$ws.Range("B401").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C401").Value2 = "Coquimbo"
$ws.Range("D401").Value2 = 44841
$ws.Range("E401").Value2 = 4
$ws.Range("F401").Value2 = 100112009
$ws.Range("G401").Value2 = "Acelga"
$ws.Range("H401").Value2 = "Sin especificar"
$ws.Range("I401").Value2 = "Primera"
$ws.Range("J401").Value2 = 2520
$ws.Range("K401").Value2 = 650
$ws.Range("L401").Value2 = 700
$ws.Range("M401").Value2 = 675
$ws.Range("N401").Value2 = "`$/atado 1,5 a 2 kilos"
$ws.Range("O401").Value2 = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P401").Value2 = 338
$ws.Range("Q401").Value2 = 2
$ws.Range("R401").Value2 = "Hortaliza"

# --- New row 402 ("Segunda") ---
$ws.Range("A402").Value2 = 8
$ws.Range("B402").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C402").Value2 = "Coquimbo"
$ws.Range("D402").Value2 = 44841
$ws.Range("E402").Value2 = 4
$ws.Range("F402").Value2 = 100112009
$ws.Range("G402").Value2 = "Acelga"
$ws.Range("H402").Value2 = "Sin especificar"
$ws.Range("I402").Value2 = "Segunda"
$ws.Range("J402").Value2 = 1400
$ws.Range("K402").Value2 = 550
$ws.Range("L402").Value2 = 600
$ws.Range("M402").Value2 = 575
$ws.Range("N402").Value2 = "`$/atado 1,5 a 2 kilos"
$ws.Range("O402").Value2 = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P402").Value2 = 288
$ws.Range("Q402").Value2 = 2
$ws.Range("R402").Value2 = "Hortaliza"

# Make sure the date cells keep/get the date number format used
# throughout column D (numFmtId 165 -> "YYYY-MM-DD HH:MM:SS", style index 2
# in the original workbook); the row insert already propagates this from
# the row above, but set it explicitly for safety.
$ws.Range("D401:D402").NumberFormat = "YYYY-MM-DD HH:MM:SS"
